$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the last paragraph in the document ("Couple other things").
# Everything new is inserted after it. We keep track of its paragraph
# *index* (not just a Range/Paragraph reference) because later
# mid-document inserts can make old references unreliable to re-derive
# positions from.
# ---------------------------------------------------------------------
$anchorIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($anchorIndex)

# ---------------------------------------------------------------------
# Helper fragments for raw-OOXML paragraph insertion (used where a
# paragraph needs more than one run, e.g. text that was typed in two
# separate pieces by the original author and so survives as two runs).
# ---------------------------------------------------------------------
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-RawParagraphAfter($afterPara, [string]$paraInnerXml) {
    # Create a placeholder paragraph right after $afterPara (it will
    # naturally inherit $afterPara's paragraph formatting / numbering).
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $placeholder = $d.Paragraphs.Item($d.Paragraphs.Count)
    $placeholderRange = $d.Range($placeholder.Range.Start, $placeholder.Range.End)

    $xml = $pkgHeader + '<w:p>' + $paraInnerXml + '</w:p>' + $pkgFooter
    $placeholderRange.InsertXML($xml)

    # InsertXML adds the new paragraph but leaves the old placeholder's
    # paragraph mark dangling as an extra empty paragraph right after
    # it; merge it away by deleting that now-redundant paragraph mark.
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $markPos = $newPara.Range.End
    $markRange = $d.Range($markPos - 1, $markPos)
    $markRange.Delete()

    return $d.Paragraphs.Item($d.Paragraphs.Count)
}

# ---------------------------------------------------------------------
# 1) "Added portal to school section" -- single-run bullet.
#    Inserted directly after "Couple other things" so it naturally
#    inherits that paragraph's ListParagraph/numPr formatting.
# ---------------------------------------------------------------------
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$bullet1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$bullet1.Range.Text = "Added portal to school section"

# ---------------------------------------------------------------------
# 2) "Added better design to " + "coding section" -- two-run bullet.
# ---------------------------------------------------------------------
$innerXml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Added better design to </w:t></w:r>' +
            '<w:r><w:t>coding section</w:t></w:r>'
$bullet2 = Insert-RawParagraphAfter $bullet1 $innerXml

# ---------------------------------------------------------------------
# 3) "Removed maininfo.css ..." -- single-run bullet.
# ---------------------------------------------------------------------
$r = $bullet2.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$bullet3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$bullet3.Range.Text = "Removed maininfo.css since the only page that required it was moved to general"

# ---------------------------------------------------------------------
# 4) "Fixed bug where main page bg ..." + " (so there was white at the
#    bottom)" -- two-run bullet.
# ---------------------------------------------------------------------
$innerXml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
            '<w:r><w:t>Fixed bug where main page bg would not scroll with content in desktop mode when scrolling was needed.</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> (so there was white at the bottom)</w:t></w:r>'
$bullet4 = Insert-RawParagraphAfter $bullet3 $innerXml

# ---------------------------------------------------------------------
# 5) "Beta 2.10" -- plain heading paragraph (no list formatting),
#    inserted between "Couple other things" and "Added portal to
#    school section" (i.e. right after $anchor, at $anchorIndex + 1).
#    Built last so that the list items above could first inherit
#    $anchor's list formatting cleanly; its own position is derived
#    from the stable $anchorIndex rather than from Paragraphs.Count,
#    since this insertion happens in the middle of the document.
# ---------------------------------------------------------------------
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$betaPara = $d.Paragraphs.Item($anchorIndex + 1)
$betaPara.Range.ListFormat.RemoveNumbers()
$betaPara.Style = "Normal"
$betaPara.Range.Text = "Beta 2.10"
